$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the empty gap rows (3:9) so that the existing entries (previously on
#    rows 10-14) shift up to rows 3-7.
$ws.Range("3:9").Delete()

# 2. Duplicate the formatting of the last existing row (now row 7) onto the two
#    brand new rows (8 and 9) so dates keep their date style and names keep
#    their style, then fill in the new contribution data.
$ws.Range("B7:E7").Copy()
$ws.Range("B8:E9").PasteSpecial(-4122)

$ws.Range("B8").Value = 44201
$ws.Range("C8").Value = "François"
$ws.Range("D8").Value = 17
$ws.Range("E8").Value = "Implémentation meilleur système"

$ws.Range("B9").Value = 44202
$ws.Range("C9").Value = "François"
$ws.Range("D9").Value = 17
$ws.Range("E9").Value = "Finalisation"

# 3. Update the contribution label of the row that used to say "Tout".
$ws.Range("E6").Value = "Implémentation système de base"

# 4. Resize column E to fit the new (longer) text.
$ws.Columns("E:E").AutoFit()

# 5. Reorder / recolor the 4 existing conditional-formatting rules on column C
#    (keep the same 4 dxf entries - just edit them in place - instead of
#    deleting/recreating them).
$rng = $ws.Range("C1:C1048576")
$fcs = $rng.FormatConditions

$fcs.Item(4).SetFirstPriority()
$fcs.Item(3).SetFirstPriority()
$fcs.Item(2).SetFirstPriority()
$fcs.Item(1).SetFirstPriority()

$fcs.Item(1).Formula1 = '"Lucas"'
$fcs.Item(1).Interior.Color = 49407

$fcs.Item(2).Formula1 = '"Thomas"'
$fcs.Item(2).Interior.Color = 16711935

$fcs.Item(3).Formula1 = '"Anthony"'
$fcs.Item(3).Interior.Color = 10498160

$fcs.Item(4).Formula1 = '"François"'
$fcs.Item(4).Interior.Color = 12611584

$fcs.Item(1).ModifyAppliesToRange($ws.Range("C1:C9"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("C1:C9"))
$fcs.Item(3).ModifyAppliesToRange($ws.Range("C1:C9"))
$fcs.Item(4).ModifyAppliesToRange($ws.Range("C1:C9"))

# 6. Leave the selection where the user ended up after typing the new rows.
$ws.Range("E11").Select()
